$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "edit2"
$ws.Range("B4").Value = "riya-morankar"
$ws.Range("C4").Value = "Merged"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2025-06-18"
$ws.Range("F4").Value = "N/A"
